$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (column C) date for every existing data row (2..468)
# from 45202 (2023-10-03) to 45203 (2023-10-04).
for ($r = 2; $r -le 468; $r++) {
    $ws.Cells.Item($r, 3).Value = 45203
}

# Row 468 picks up an explicit row height (matches the rest of the sheet).
$ws.Rows.Item(468).RowHeight = 15

# Append the new record as row 469.
$ws.Cells.Item(469, 1).Value = "A 47216-2023"

$ws.Cells.Item(469, 2).Value = 45202
$ws.Cells.Item(469, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(469, 3).Value = 45203
$ws.Cells.Item(469, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(469, 4).Value = "JÖNKÖPINGS LÄN"
$ws.Cells.Item(469, 5).Value = "VÄRNAMO"

$ws.Cells.Item(469, 7).Value = 0.9
$ws.Cells.Item(469, 8).Value = 0
$ws.Cells.Item(469, 9).Value = 0
$ws.Cells.Item(469, 10).Value = 0
$ws.Cells.Item(469, 11).Value = 0
$ws.Cells.Item(469, 12).Value = 0
$ws.Cells.Item(469, 13).Value = 0
$ws.Cells.Item(469, 14).Value = 0
$ws.Cells.Item(469, 15).Value = 0
$ws.Cells.Item(469, 16).Value = 0
$ws.Cells.Item(469, 17).Value = 0

$ws.Cells.Item(469, 18).Value = ""
$ws.Cells.Item(469, 18).WrapText = $true
